$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the "Repollo" series for
# Terminal Hortofrutícola Agro Chillán. Insert a fresh row at row 97 —
# this pushes the existing rows 97:178 down to 98:179 (dimension grows
# to A1:R179), matching the rest of the diff, which is just every row
# from 97 on shifted down by one.
$ws.Rows(97).Insert()

# Populate the newly inserted row 97 with the new record. Every column
# other than the price/volume/date fields repeats the same
# Terminal Hortofrutícola Agro Chillán / Ñuble / Repollo / Crespo record /
# Primera / Provincia de Diguillín template used throughout this block.
$ws.Cells.Item(97, 1).Value = 7
$ws.Cells.Item(97, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(97, 3).Value = "Ñuble"
$ws.Cells.Item(97, 4).Value = 44589
$ws.Cells.Item(97, 5).Value = 16
$ws.Cells.Item(97, 6).Value = 100112006
$ws.Cells.Item(97, 7).Value = "Repollo"
$ws.Cells.Item(97, 8).Value = "Crespo record"
$ws.Cells.Item(97, 9).Value = "Primera"
$ws.Cells.Item(97, 10).Value = 300
$ws.Cells.Item(97, 11).Value = 650
$ws.Cells.Item(97, 12).Value = 700
$ws.Cells.Item(97, 13).Value = 675
$ws.Cells.Item(97, 14).Value = "$/unidad"
$ws.Cells.Item(97, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(97, 16).Value = 675
$ws.Cells.Item(97, 17).Value = 1
$ws.Cells.Item(97, 18).Value = "Hortaliza"
